# Fix: added rolls in att-bulk-format
#
# The attendance bulk-upload template originally had columns:
#   A=Student Name, B=Student ID, C..G = attendance dates
# This change inserts two new columns after "Student ID" for
# "Class Roll" and "Univ Roll", pushing the date columns from C:G to E:I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at C:D - this shifts the existing date columns
# (previously C:G) two places to the right (now E:I), carrying their
# values, number formats and styles with them.
$ws.Range("C1:D1").EntireColumn.Insert()

# Populate the two new header cells.
$ws.Range("C1").Value = "Class Roll"
$ws.Range("D1").Value = "Univ Roll"

# Restore the selection to where the author ended up after editing.
$ws.Range("I8").Select()
